$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 118.55556
$ws.Range("I2").Value = 61
$ws.Range("J2").Value = 190.5
$ws.Range("K2").Value = 61
$ws.Range("L2").Value = 190.5
$ws.Range("M2").Value = 52
$ws.Range("N2").Value = -416.5

$ws.Range("H40").Value = 1443.5
$ws.Range("J40").Value = 1615.2
$ws.Range("L40").Value = 1615.2
$ws.Range("N40").Value = -1965.2

$ws.Range("H132").Value = 1731.3116
$ws.Range("I132").Value = 1346.3442
$ws.Range("K132").Value = 4039.0326
$ws.Range("M132").Value = -1509.0326

$ws.Range("H140").Value = 55450
$ws.Range("J140").Value = 55450
$ws.Range("L140").Value = 55450
$ws.Range("N140").Value = -65810

$ws.Range("H141").Value = 5272.775
$ws.Range("I141").Value = 2841.6487
$ws.Range("K141").Value = 8524.946100000001
$ws.Range("M141").Value = -3344.946100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19901.457
$ws.Range("I32").Value = 19146.977
$ws.Range("J32").Value = 36500
$ws.Range("K32").Value = 19146.977
$ws.Range("L32").Value = 36500
$ws.Range("M32").Value = -18859.977
$ws.Range("N32").Value = -37074

$ws.Range("H132").Value = 2803.4426
$ws.Range("I132").Value = 2793.647
$ws.Range("J132").Value = 2853.4
$ws.Range("K132").Value = 8380.940999999999
$ws.Range("L132").Value = 8560.200000000001
$ws.Range("M132").Value = -5850.940999999999
$ws.Range("N132").Value = -13620.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 9800
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 9800
$ws.Range("K6").Value = 0
$ws.Range("L6").ClearContents()
$ws.Range("M6").Value = 9800
$ws.Range("N6").Value = -10026

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2090.25
$ws.Range("I16").Value = 2453.6667
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 2453.6667
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -2166.6667
$ws.Range("N16").Value = -1574

$ws.Range("H31").Value = 30309944
$ws.Range("I31").Value = 83334970
$ws.Range("K31").Value = 83334970
$ws.Range("M31").Value = -83334675

$ws.Range("H34").Value = 30309944
$ws.Range("I34").Value = 83334970
$ws.Range("K34").Value = 83334970
$ws.Range("M34").Value = -83334768

$ws.Range("H58").Value = 40001364
$ws.Range("I58").Value = 71429550
$ws.Range("J58").Value = 1854.3636
$ws.Range("K58").Value = 71429550
$ws.Range("L58").Value = 1854.3636
$ws.Range("M58").Value = -71429347
$ws.Range("N58").Value = -2260.3636

$ws.Range("H113").Value = 2090.25
$ws.Range("I113").Value = 2453.6667
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2453.6667
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -283.6667000000002
$ws.Range("N113").Value = -5340

$ws.Range("H134").Value = 2301.4
$ws.Range("I134").Value = 2209.647
$ws.Range("K134").Value = 6628.941
$ws.Range("M134").Value = -4093.941

$ws.Range("H136").Value = 40001364
$ws.Range("I136").Value = 71429550
$ws.Range("J136").Value = 1854.3636
$ws.Range("K136").Value = 214288650
$ws.Range("L136").Value = 5563.0908
$ws.Range("M136").Value = -214286100
$ws.Range("N136").Value = -10663.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 698.2045000000001
$ws.Range("I5").Value = 469.08572
$ws.Range("J5").Value = 1589.2222
$ws.Range("K5").Value = 1407.25716
$ws.Range("L5").Value = 4767.6666
$ws.Range("M5").Value = -1295.25716
$ws.Range("N5").Value = -4991.6666

$ws.Range("H116").Value = 1552.8
$ws.Range("I116").Value = 2071.3333
$ws.Range("J116").Value = 775
$ws.Range("K116").Value = 6213.999899999999
$ws.Range("L116").Value = 2325
$ws.Range("M116").Value = -2771.999899999999
$ws.Range("N116").Value = -9209

$ws.Range("H122").Value = 1130.85
$ws.Range("I122").Value = 593.3333
$ws.Range("J122").Value = 1937.125
$ws.Range("K122").Value = 5339.9997
$ws.Range("L122").Value = 17434.125
$ws.Range("M122").Value = -2889.9997
$ws.Range("N122").Value = -22334.125

$ws.Range("H135").Value = 698.2045000000001
$ws.Range("I135").Value = 469.08572
$ws.Range("J135").Value = 1589.2222
$ws.Range("K135").Value = 4221.771479999999
$ws.Range("L135").Value = 14302.9998
$ws.Range("M135").Value = -1686.771479999999
$ws.Range("N135").Value = -19372.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 48000
$ws.Range("J32").Value = 48000
$ws.Range("L32").Value = 48000
$ws.Range("N32").Value = -48592

$ws.Range("H42").Value = 40799
$ws.Range("J42").Value = 40799
$ws.Range("L42").Value = 40799
$ws.Range("N42").Value = -41769

$ws.Range("H70").Value = 4975.5557
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 4968.5713
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 4968.5713
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = -5508.5713

$ws.Range("H73").Value = 4975.5557
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 4968.5713
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 4968.5713
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = -6840.5713

$ws.Range("H115").Value = 40799
$ws.Range("J115").Value = 40799
$ws.Range("L115").Value = 40799
$ws.Range("N115").Value = -43149

$ws.Range("H138").Value = 49429
$ws.Range("J138").Value = 49429
$ws.Range("L138").Value = 49429
$ws.Range("N138").Value = -59709

$ws.Range("H139").Value = 10000
$ws.Range("I139").Value = 10000
$ws.Range("K139").Value = 10000
$ws.Range("M139").Value = -4860

$ws.Range("H140").Value = 44400
$ws.Range("J140").Value = 44400
$ws.Range("L140").Value = 44400
$ws.Range("N140").Value = -54760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1232.56
$ws.Range("I61").Value = 962.7
$ws.Range("J61").Value = 2312
$ws.Range("K61").Value = 962.7
$ws.Range("L61").Value = 2312
$ws.Range("M61").Value = -760.7
$ws.Range("N61").Value = -2716

$ws.Range("H113").Value = 1232.56
$ws.Range("I113").Value = 962.7
$ws.Range("J113").Value = 2312
$ws.Range("K113").Value = 962.7
$ws.Range("L113").Value = 2312
$ws.Range("M113").Value = 1207.3
$ws.Range("N113").Value = -6652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3798.0266
$ws.Range("I136").Value = 4589.4614
$ws.Range("J136").Value = 2008.6957
$ws.Range("K136").Value = 13768.3842
$ws.Range("L136").Value = 6026.0871
$ws.Range("M136").Value = -11218.3842
$ws.Range("N136").Value = -11126.0871
